$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101, shifting existing rows 101:120 down to 102:121
$ws.Rows("101:101").Insert()

# Populate the newly inserted row 101 with the new weekly record
$ws.Range("A101").Value = 10
$ws.Range("B101").Value = "Vega Modelo de Temuco"
$ws.Range("C101").Value = "La Araucanía"
$ws.Range("D101").Value = 45015
$ws.Range("E101").Value = 9
$ws.Range("F101").Value = 100112030
$ws.Range("G101").Value = "Poroto granado"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 40
$ws.Range("K101").Value = 42000
$ws.Range("L101").Value = 42000
$ws.Range("M101").Value = 42000
$ws.Range("N101").Value = "$/saco 25 kilos"
$ws.Range("O101").Value = "Región de La Araucanía"
$ws.Range("P101").Value = 1680
$ws.Range("Q101").Value = 25
$ws.Range("R101").Value = "Hortaliza"
